$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 8 new blank rows starting at row 3; this pushes the existing
# row 2 ("Oct 7 2020" vs Chennai Super Kings) down to row 4.
$ws.Range("A3:A10").EntireRow.Insert()

# Final table (rows 2-10), all columns A..K.
$data = @(
    @(" Oct 3 2020", " Sharjah", "Capitals won by 18 runs", "Kolkata Knight Riders", "Delhi Capitals", "Sunil Narine ", "3", "5", "0", "0", "60.00"),
    @(" Oct 24 2020", " Abu Dhabi", "KKR won by 59 runs", "Kolkata Knight Riders", "Delhi Capitals", "Sunil Narine ", "64", "32", "6", "4", "200.00"),
    @(" Oct 7 2020", " Abu Dhabi", "KKR won by 10 runs", "Kolkata Knight Riders", "Chennai Super Kings", "Sunil Narine ", "17", "9", "1", "1", "188.88"),
    @(" Oct 29 2020", " Dubai (DSC)", "Super Kings won by 6 wickets", "Kolkata Knight Riders", "Chennai Super Kings", "Sunil Narine ", "7", "7", "0", "1", "100.00"),
    @(" Oct 26 2020", " Sharjah", "Kings XI won by 8 wickets (with 7 balls remaining)", "Kolkata Knight Riders", "Kings XI Punjab", "Sunil Narine ", "6", "4", "1", "0", "150.00"),
    @(" Sep 26 2020", " Abu Dhabi", "KKR won by 7 wickets (with 12 balls remaining)", "Kolkata Knight Riders", "Sunrisers Hyderabad", "Sunil Narine ", "0", "2", "0", "0", "0.00"),
    @(" Nov 1 2020", " Dubai (DSC)", "KKR won by 60 runs", "Kolkata Knight Riders", "Rajasthan Royals", "Sunil Narine ", "0", "2", "0", "0", "0.00"),
    @(" Sep 30 2020", " Dubai (DSC)", "KKR won by 37 runs", "Kolkata Knight Riders", "Rajasthan Royals", "Sunil Narine ", "15", "14", "2", "1", "107.14"),
    @(" Sep 23 2020", " Abu Dhabi", "Mumbai won by 49 runs", "Kolkata Knight Riders", "Mumbai Indians", "Sunil Narine ", "9", "10", "0", "1", "90.00")
)

# Write every value as a quoted-string formula so Excel keeps the cell as
# text (matching the source data where numeric-looking values such as
# "17" or "188.88" are stored as text, not numbers).
for ($i = 0; $i -lt $data.Length; $i++) {
    $rowNum = $i + 2
    $values = $data[$i]
    for ($j = 0; $j -lt $values.Length; $j++) {
        $colNum = $j + 1
        $escaped = $values[$j].Replace('"', '""')
        $ws.Cells.Item($rowNum, $colNum).Formula = '="' + $escaped + '"'
    }
}

# Flatten every formula in the table down to a literal value in one pass
# so the saved file has plain text cells (no formulas, no new styles).
$fullRange = $ws.Range("A1:K10")
$fullRange.Copy()
$fullRange.PasteSpecial(-4163)
$excel.CutCopyMode = $false
